$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 43; this shifts existing rows 43:134
# down to 44:135 and extends the used range to A1:T135.
$ws.Rows.Item(43).Insert()

# Populate the newly inserted row 43 with the new weekly record.
$ws.Cells.Item(43, 1).Value = 10
$ws.Cells.Item(43, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(43, 3).Value = "La Araucanía"
$ws.Cells.Item(43, 4).Value = 45274
$ws.Cells.Item(43, 5).Value = 9
$ws.Cells.Item(43, 6).Value = "Fruta"
$ws.Cells.Item(43, 7).Value = 100108
$ws.Cells.Item(43, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(43, 9).Value = 100108007
$ws.Cells.Item(43, 10).Value = "Coco"
$ws.Cells.Item(43, 11).Value = "Sin especificar"
$ws.Cells.Item(43, 12).Value = "Primera"
$ws.Cells.Item(43, 13).Value = 50
$ws.Cells.Item(43, 14).Value = 32000
$ws.Cells.Item(43, 15).Value = 32000
$ws.Cells.Item(43, 16).Value = 32000
$ws.Cells.Item(43, 17).Value = "`$/malla 20 unidades"
$ws.Cells.Item(43, 18).Value = "Perú"
$ws.Cells.Item(43, 19).Value = 1600
$ws.Cells.Item(43, 20).Value = 20
